$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Basic" row (row 4) values
$ws.Range("D4").Value = 13
$ws.Range("E4").Value = 5

# Move the active selection to E5 (cursor moved after editing E4)
$ws.Range("E5").Select()
